# Generate Report for Handback
# The handback process detected that the handback file name did not match
# the handoff file name for the "1f6ca00a-ebf9-4deb-afcf-d72b7d2510cb.md"
# entry, so the status is updated to reflect the failure and an error
# detail message is recorded for each locale sheet.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$failMessage = "Handback transform failed"

# Overview sheet: row 3 corresponds to the 1f6ca00a-ebf9-4deb-afcf-d72b7d2510cb.md file.
# Both the zh-cn (B3) and de-de (C3) status columns move from "Ready for handoff"
# to "Handback transform failed".
$overview.Range("B3").Value = $failMessage
$overview.Range("C3").Value = $failMessage

# The locale sheets (zh-cn / de-de) also show the status for the same file in
# their own "Status" column (C3), which shared the same string value.
$zhcn.Range("C3").Value = $failMessage
$dede.Range("C3").Value = $failMessage

# zh-cn sheet: row 3 (same file) gets an Error Detail entry in column K.
$zhcn.Range("K3").Value = "Handback file name: co3xeyno.z5z is different with handoff file name: 1f6ca00a-ebf9-4deb-afcf-d72b7d2510cb.405672d2857fde8453891dafe7307e5ac8fa54f8.zh-cn."

# de-de sheet: row 3 (same file) gets an Error Detail entry in column K.
$dede.Range("K3").Value = "Handback file name: co3xeyno.z5z is different with handoff file name: 1f6ca00a-ebf9-4deb-afcf-d72b7d2510cb.405672d2857fde8453891dafe7307e5ac8fa54f8.de-de."
